# edit.ps1 -- applies the "cosmos" -> "chemistry" rewrite described by the
# supplied unified diff to the active document.
#
# Strategy:
#  * Plain in-run text swaps are done with Find.Execute (wdReplaceOne),
#    which keeps the original run's rPr untouched.
#  * Where the diff *adds* brand new runs (a lone "." run followed by a new
#    sentence run) we insert the text right after the just-replaced range
#    and then force a run boundary (an immaterial Bold on/off toggle) before
#    stamping the run's font explicitly, so the produced OOXML has its own
#    <w:rPr> identical to its neighbours instead of silently merging back
#    into the preceding run.

function New-RunAfter {
    param($Doc, $Pos, $Text, $FontName, $FontSize, $HasSize, $FontColor)
    $ip = $Doc.Range($Pos, $Pos)
    $ip.InsertAfter($Text)
    # Force a distinct run: toggle a property on then straight back off.
    $ip.Font.Bold = 1
    $ip.Font.Bold = 0
    $ip.Font.Name = $FontName
    if ($HasSize -eq 1) {
        $ip.Font.Size = $FontSize
    }
    $ip.Font.Color = $FontColor
    return $Pos + $Text.Length
}

function Replace-Once {
    param($Doc, $Find, $Replace)
    $rng = $Doc.Content
    $ok = $rng.Find.Execute($Find, $true, $false, $false, $false, $false, $true, 1, $false, $Replace, 2)
    return $rng
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------
Replace-Once $d "Unraveling the Mysteries of the Cosmos" "Unveiling the Enigmatic World of Chemistry" | Out-Null

# ---------------------------------------------------------------------
# Author name: "Aria Dimitrov" -> "Dr" + "." + " Caroline Stevens"
# ---------------------------------------------------------------------
$rng = Replace-Once $d "Aria Dimitrov" "Dr"
$pos = $rng.End
$pos = New-RunAfter $d $pos "." "Calibri" 18 1 0
$pos = New-RunAfter $d $pos " Caroline Stevens" "Calibri" 18 1 0

# ---------------------------------------------------------------------
# Email: "cosmos_quest@spaceworld" -> "caroline" ; "." unchanged ;
# "net" -> "stevens@schoolmail" + "." + "com"
# ---------------------------------------------------------------------
Replace-Once $d "cosmos_quest@spaceworld" "caroline" | Out-Null
$rng = Replace-Once $d "net" "stevens@schoolmail"
$pos = $rng.End
$pos = New-RunAfter $d $pos "." "Calibri" 16 1 0
$pos = New-RunAfter $d $pos "com" "Calibri" 16 1 0

# ---------------------------------------------------------------------
# Main body paragraph (size 24 half-points == 12pt)
# ---------------------------------------------------------------------
Replace-Once $d "The cosmos, an enigmatic tapestry of celestial wonders, has captivated humanity since the dawn of time" "Chemistry, the study of matter, offers a path to unravel the hidden mysteries of the physical world" | Out-Null

Replace-Once $d " From ancient astronomers gazing upon the star-studded night sky to modern scientists exploring the outer reaches of the universe, our quest to unravel its mysteries has been an enduring and awe-inspiring pursuit" " It's a science that uncovers the intricacies of tiny particles, revealing an intriguing ballet of atoms and molecules" | Out-Null

$rng = Replace-Once $d " The cosmos, vast and mysterious, serves as an infinite canvas upon which the grand spectacle of cosmic events unfolds, beckoning us to ponder upon our place in this intricate cosmic symphony" " With each experiment, we peek behind the curtain of the everyday, witnessing the magic of chemical reactions and learning the secrets of substances all around us"
$pos = $rng.End
$pos = New-RunAfter $d $pos "." "Calibri" 12 1 0
$pos = New-RunAfter $d $pos " Our lives are woven with countless feats of chemistry -- from the food we eat to the medicines that heal us" "Calibri" 12 1 0

Replace-Once $d "As we traverse this celestial odyssey, we encounter cosmic phenomena that defy our understanding" "Chemistry allows us to explore the diverse tapestry of materials, from the stardust of diamond to the elasticity of rubber" | Out-Null

Replace-Once $d " Supermassive black holes, enigmatic entities lurking at the heart of galaxies, possess a gravitational pull so intense that not even light can escape their clutches" " It's a journey through substances and their interactions, a chronicle of transformations and creations" | Out-Null

Replace-Once $d " Neurons, the intricate building blocks of human consciousness, orchestrate a symphony of electrical impulses, enabling us to perceive and navigate the world around us" " Each element tells a story, and each reaction is a chapter in the epic saga of chemistry's symphony" | Out-Null

Replace-Once $d " Quantum mechanics, the perplexing realm of the infinitely small, presents us with paradoxes that challenge our very perception of reality" " To understand chemistry is to decipher the language of the physical world, unveiling the dance of particles that defines our existence" | Out-Null

Replace-Once $d "These cosmic mysteries, both grand and intricate, ignite a burning curiosity within us" "Chemistry's enigma lies in its universality" | Out-Null

Replace-Once $d " They prompt us to delve deeper, to seek answers to questions that have perplexed humanity for millennia" " It transcends species and generations, binding all life together in a web of chemical connections" | Out-Null

Replace-Once $d " Our relentless pursuit of understanding has led to profound discoveries, expanding our knowledge of the universe and our place within it" " From the vibrant hues of nature to the intricate machinery of cells, chemistry is a shared language that manifests in countless forms" | Out-Null

$rng = Replace-Once $d " Yet, as we unravel one enigma, another emerges, beckoning us to continue our exploration, forever captivated by the boundless wonders of the cosmos" " Its universality challenges us to seek patterns in diversity, revealing the interconnectedness of all things"
$pos = $rng.End
$pos = New-RunAfter $d $pos "." "Calibri" 12 1 0
$pos = New-RunAfter $d $pos " In the realm of chemistry, the ordinary becomes extraordinary as we witness the mundane transformed into a universe of atoms and molecules" "Calibri" 12 1 0

# ---------------------------------------------------------------------
# Summary paragraph (default size -- no explicit w:sz)
# ---------------------------------------------------------------------
Replace-Once $d "Our journey into the cosmos, driven by an insatiable curiosity, has unveiled cosmic wonders that defy comprehension" "Chemistry, the study of matter, unveils the captivating world of substances and reactions" | Out-Null

Replace-Once $d " From black holes devouring matter and energy to neurons facilitating the marvels of consciousness, the universe presents us with perplexing phenomena that challenge our understanding" " It takes us on a journey through the diverse tapestry of materials, allowing us to decipher the language of the physical world" | Out-Null

Replace-Once $d " Despite the challenges, our unwavering pursuit of knowledge has led to groundbreaking discoveries, propelling us forward in our quest to unravel the intricate mysteries of the cosmos" " Through chemistry, we explore the enigmatic interactions of atoms and molecules, witnessing the epic saga of transformations and creations" | Out-Null

$rng = Replace-Once $d " Our exploration will continue, forever fueled by the allure of the unknown, as we strive to comprehend the enigmatic tapestry of the universe that envelops us" " Its universality binds all life together, inviting us to uncover the interconnectedness of all things"
$pos = $rng.End
$pos = New-RunAfter $d $pos "." "Calibri" 0 0 0
$pos = New-RunAfter $d $pos " With every experiment, chemistry invites us to question, explore, and discover the hidden marvels of our physical world" "Calibri" 0 0 0

# ---------------------------------------------------------------------
# Trailing empty paragraph
# ---------------------------------------------------------------------
$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRng.InsertParagraphAfter()

Write-Output "done"
